# Add SF_PLAY_MODE to test data.
#
# The original edit inserted a new column immediately to the left of the
# existing "AV_FILE" column (old column AV), shifting every column from
# AV onward one place to the right (old AV..AZ -> new AW..BA). The new
# column is used for a single header cell "SF_PLAY_MODE" in row 1; the
# data rows (2-4) are left blank in the new column, same as Excel does
# when inserting a column (new cells inherit the left neighbour's style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before AV (old AV_FILE column), pushing AV:AZ to AW:BA.
$ws.Range("AV1").EntireColumn.Insert()

# New header cell for the inserted column.
$ws.Range("AV1").Value = "SF_PLAY_MODE"

# Give the new column a sensible explicit width (close to its neighbours,
# the engine quantizes ColumnWidth to 1/6-character steps on save, so this
# is the closest input to the original column's ~17.9-character width).
$ws.Range("AV1").ColumnWidth = 17

# The worksheet's hidden _FilterDatabase name covers the header/data rows;
# it needs to grow by one column (AY4 -> AZ4) to keep tracking the full
# A1:<lastcol>4 range now that a column was inserted inside it.
$names = $wb.Names
for ($i = 1; $i -le $names.Count(); $i++) {
    $n = $names.Item($i)
    if ($n.Name() -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$AZ`$4"
    }
}
